$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new timesheet entries (rows 42 and 43) -----------------------
# Row 42 mirrors the formatting of row 37 (style index 6 pattern),
# row 43 mirrors the formatting of row 41 (style index 4 pattern) - this
# matches the alternating cell styles already present throughout the sheet.

$ws.Range("C37:G37").Copy()
$ws.Range("C42:G42").PasteSpecial(-4122)

$ws.Range("C41:G41").Copy()
$ws.Range("C43:G43").PasteSpecial(-4122)

# Row 42: SlNo 36, 29-Aug-2021, Java Traning / Learing java 8 Features, 5 hrs
$ws.Range("C42").Value = 36
$ws.Range("D42").Value = 44437
$ws.Range("E42").Value = "Java Traning "
$ws.Range("F42").Value = "Learing java 8 Features"
$ws.Range("G42").Value = 5

# Row 43: SlNo 37, 30-Aug-2021, Java / File Handing, 6 hrs
# (set F43 before E43 so the new shared strings are appended in the same
# order as the source workbook: "Java Traning ", "File Handing ", "Java ")
$ws.Range("C43").Value = 37
$ws.Range("D43").Value = 44438
$ws.Range("F43").Value = "File Handing "
$ws.Range("E43").Value = "Java "
$ws.Range("G43").Value = 6

# --- Update the view: scroll down and select K44 --------------------------
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K44").Select()
